$d = $word.ActiveDocument

# 1. Replace the document title text.
$d.Content.Find.Execute(
    "2.2 - Debate I", $true, $false, $false, $false, $false,
    $true, 1, $false, "Placeholder - Check Back Later", 2
) | Out-Null

# 2. Remove the trailing " :::" that follows "...general edification later."
#    (deletes the two runs holding the space and the ":::" marker).
$d.Content.Find.Execute(
    " :::", $true, $false, $false, $false, $false,
    $true, 1, $false, "", 2
) | Out-Null
